$d = $word.ActiveDocument

# Locate the redundant word "труб " ("pipes ") that needs to be removed from the
# sentence "...из %PM% труб – %PL% км." so it reads "...из %PM% – %PL% км."
$rng = $d.Content
$rng.Find.Execute("труб ", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End

# Remove the word (and its trailing space) in place.
$delRng = $d.Range($start, $end)
$delRng.Delete()

# Mark the edit location with a "_GoBack" bookmark, as Word does automatically
# when a document is edited and saved.
$bmRng = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $bmRng)
